$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 6788.088191908992
$ws.Range("G3").Value = 17086.97919129642
$ws.Range("G4").Value = 22231.83681584428
$ws.Range("G5").Value = 6198.711515966546
$ws.Range("G6").Value = 0.4978133868217219
$ws.Range("G7").Value = 0.5797070920391136
$ws.Range("F8").Value = 0.5709397031663825
$ws.Range("G8").Value = 0.5545366382550266
$ws.Range("G9").Value = 0.4927461019866803
$ws.Range("G10").Value = 0.4218885172214436
$ws.Range("G11").Value = 0.4965821404509979
$ws.Range("G12").Value = 0.5480300352097685
$ws.Range("G13").Value = 0.5123615864059685
